$wb = $excel.ActiveWorkbook

# This workbook is refreshed from a live "next bus arrival" data source
# (commit: "Sync file from Google Drive"). NextBus1 is unchanged; NextBus2
# and NextBus3 each get refreshed ETA / minutes-to-arrival figures, one
# fewer upcoming bus is reported (last row is dropped), and what used to
# be row 7's bus now shows up as row 6 with its own refreshed figures.

foreach ($sheetName in @("NextBus2", "NextBus3")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Refreshed EstimatedTimeOfArrival (col F) and MinutesToArrival (col O)
    # for the buses that stay in rows 2-5.
    $ws.Cells.Item(2, 6).Value = 45688.61592592593
    $ws.Cells.Item(2, 15).Value = 26

    $ws.Cells.Item(3, 6).Value = 45688.61098379629
    $ws.Cells.Item(3, 15).Value = 19

    $ws.Cells.Item(4, 6).Value = 45688.62466435185
    $ws.Cells.Item(4, 15).Value = 38

    $ws.Cells.Item(5, 6).Value = 45688.61802083333
    $ws.Cells.Item(5, 15).Value = 29

    # Row 6 now carries what used to be row 7's bus (BusNo, ETA,
    # OriginCode, MinutesToArrival refreshed); destination columns (C/D)
    # stay as they were.
    $ws.Cells.Item(6, 2).Value = 74
    $ws.Cells.Item(6, 6).Value = 45688.61518518518
    $ws.Cells.Item(6, 11).Value = 11379
    $ws.Cells.Item(6, 15).Value = 25

    # The old row 7 is gone entirely now.
    $ws.Rows(7).Delete()
}
